$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AV (column 48). This shifts the existing
# "nom" (AV) and "url_produit" (AW) columns one position to the right,
# landing on AW and AX respectively, and preserves their styles/values.
$ws.Columns("AV").Insert()

# New header timestamp for the freshly inserted price-history column.
$ws.Range("AV1").Value = "2026-01-29 19:23:33"

# For every product row that still had a price tracked in the last
# existing history column (AU), carry that same last-known price
# forward into the newly inserted AV column (rows 2-80). Rows whose AU
# cell was already empty (no longer tracked) are left untouched, which
# keeps the newly inserted AV cell empty for them as well (rows 81-206).
for ($row = 2; $row -le 80; $row++) {
    $lastPrice = $ws.Cells.Item($row, 47).Value2
    $ws.Cells.Item($row, 48).Value2 = $lastPrice
}
